$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename column header "patient_id" -> "subject_id"
$ws.Range("D1").Value = "subject_id"

# Update the active cell selection to D2
$ws.Range("D2").Select()
